# Anjana's diary entry for week 7 (#407)
# Fills in the three blank diary rows (51, 53, 55) covering the week of
# 2020-02-20 (mid-term results / architecture), 2020-02-23 and 2020-02-24
# (h2 as-described vs as-implemented architecture homework).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an already-formatted diary row as a style template (date style in A,
# wrap-text entry style in B:F, mood style in G) so the new rows pick up
# the exact same cell formatting as every other entry in the table.
$ws.Range("A45:G45").Copy($ws.Range("A51:G51"))
$ws.Range("A45:G45").Copy($ws.Range("A53:G53"))
$ws.Range("A45:G45").Copy($ws.Range("A55:G55"))

# --- Dates / reused Time & Participants values (existing shared strings) ---------
$ws.Range("A51").Value = 43881
$ws.Range("B51").Value = "5pm - 7.50 pm"
$ws.Range("C51").Value = "N/A"

$ws.Range("A53").Value = 43884
$ws.Range("C53").Value = "Vaishakhi, Aman"

$ws.Range("A55").Value = 43885
$ws.Range("C55").Value = "Vaishakhi, Aman"

# --- New text entered in authoring order (keeps shared-string table append order
#     identical to the original edit) ---------------------------------------------
$ws.Range("D51").Value = "Was expecting mid-term results, looking forward to know more about the system architecture"
$ws.Range("E51").Value = "Learned three new Key Expert practices.`nLearned about the architecture of the system and the process of understanding the architecture from source code."
$ws.Range("F51").Value = "If there's no proper documented architecture, the first step in understanding the architecture from the source code is by looking at the folder/package  structure and trying to grouping various related classes together. We can start with the UML diagram and slowly try to abstract up. `nPull requests can also be a useful tool as lot of design decisions could be present in pull requests which will help us to understand the rationale of the developer and why certain things are written in a certain way. "
$ws.Range("G51").Value = "This week's homework is pretty vast and have to start working on it soon."

$ws.Range("E53").Value = "Understood both the as-described and as-implemented architecture of the system."
$ws.Range("F53").Value = "Realized the importance of having a documented version of architecture as it helps developers in the future `nThere are only few variations in the as-described and as-implemented architectures of h2, which means the h2 community has done a really good job at code reviews and maintaining the standards`n"
$ws.Range("G53").Value = "Proud to complete the most challenging part of this week's assignment"
$ws.Range("D53").Value = "Understand the architecture of the system and document the same"

$ws.Range("D55").Value = "Finish the remaining sections of the homework like pull requests, issues, state of the system etc. "
$ws.Range("E55").Value = "We divided the remaining sections of the homework and worked on it individually. Later on we collated all our findings in the report."

$ws.Range("B53").Value = "1pm - 6pm"

$ws.Range("F55").Value = "Looking for the social context was not difficult as h2 maintains good documentation in their website`nUnderstood the importance of having proper comments and explanation in the pull requests as it helped in understanding the decisions made by the contributors and the rationale behind each change"
$ws.Range("G55").Value = "Happy to complete the homework early"

$ws.Range("B55").Value = "9pm - 12am"

# --- Row heights (auto-grown by Excel when the wrapped text was entered) ---------
$ws.Rows.Item(51).RowHeight = 249.6
$ws.Rows.Item(53).RowHeight = 156
$ws.Rows.Item(55).RowHeight = 156

# --- Restore the view: scrolled down to the newly added rows, selection on
#     the last entered cell, zoomed to 100%. -------------------------------------
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.Zoom = 100
$ws.Range("F55").Select()
